# Added Linguistic Object section + prev accession number fix
#
# This script updates the "Objects" sheet of the Linked Art model-tracking
# workbook:
#   - fixes the "Previous Accession Number" and "Alternate Titles" URI
#     patterns to use the "[sequence #]" placeholder style
#   - fills in the "Series Title" / "Portfolio Title" / "Title Notes" rows
#     (which previously only had a label) with their EMu field, URI pattern
#     and notes, matching the "Alternate Titles" row above them
#   - highlights the two new/updated rows the same way other "needs
#     attention" cells are highlighted elsewhere in the sheet
#   - updates the frozen-pane scroll position / selection to point at the
#     newly edited rows
#   - widens the URI column (D) to match the Notes column (E)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects")

# --- 1. Previous Accession Number: URI pattern now carries a sequence # ---
$ws.Range("D5").Value = "object/[irn]/old-accession-number/[sequence #]"

# --- 2. Alternate Titles: URI pattern now carries a sequence # ---
$ws.Range("D9").Value = "object/[irn]/alt-title-[sequence #]"

# --- 3. Series Title row (10): fill in URI + notes, highlight the label ---
$ws.Range("D10").Value = "object/[irn]/series-title"
$ws.Range("F10").Value = "See GitHub Issues"

# --- 4. Portfolio Title row (11): fill in Mapped?, URI, AAT URI + notes ---
$ws.Range("B11").Value = "X"
$ws.Range("D11").Value = "object/[irn]/portfolio-title"
$ws.Range("E11").Value = "http://vocab.getty.edu/aat/300417225"
$ws.Range("F11").Value = "See GitHub Issues"

# --- 5. Title Notes row (12): fill in Mapped?, URI; URI value now spans 2
#        AAT links, styled with wrapped text like other multi-line notes ---
$ws.Range("B12").Value = "X"
$ws.Range("D12").Value = "object/[irn]/title-statement"
$ws.Range("C13").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E12").Value = "http://vocab.getty.edu/aat/300417212" + [char]10 + "http://vocab.getty.edu/aat/300418049"
$ws.Rows.Item(12).RowHeight = 30

# --- Highlight the "Series Title" / "Portfolio Title" labels the same way
#     the "Color"/"Meaning" key-sheet highlight cell is shaded, and apply
#     the same highlight to the new Portfolio Title AAT URI cell ---
$keyws = $wb.Worksheets.Item("Key")
$keyws.Range("A3").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$keyws.Range("A3").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$keyws.Range("A3").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = "http://vocab.getty.edu/aat/300417225"

# --- Widen the URI column (D) to match the Notes column (E) ---
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# --- Move the frozen-pane scroll position / selection to the edited rows ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A10:A11").Select()
